$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 763.5
$ws.Range("I55").Value = 55
$ws.Range("J55").Value = 999.6667
$ws.Range("K55").Value = 55
$ws.Range("L55").Value = 999.6667
$ws.Range("M55").Value = 159
$ws.Range("N55").Value = -1427.6667

$ws.Range("H106").Value = 5817.8213
$ws.Range("I106").Value = 5996.115
$ws.Range("J106").Value = 3500
$ws.Range("K106").Value = 5996.115
$ws.Range("L106").Value = 3500
$ws.Range("M106").Value = -5365.115
$ws.Range("N106").Value = -4762

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws.Range("H138").Value = 436653.3
$ws.Range("J138").Value = 513463.2
$ws.Range("L138").Value = 1540389.6
$ws.Range("N138").Value = -1550669.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1093.2
$ws.Range("I45").Value = 1028.4286
$ws.Range("K45").Value = 1028.4286
$ws.Range("M45").Value = -651.4286

$ws.Range("H132").Value = 2985.2273
$ws.Range("I132").Value = 2510.8667
$ws.Range("J132").Value = 4001.7144
$ws.Range("K132").Value = 7532.6001
$ws.Range("L132").Value = 12005.1432
$ws.Range("M132").Value = -5002.6001
$ws.Range("N132").Value = -17065.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 45905296
$ws.Range("I105").Value = 48091172
$ws.Range("K105").Value = 48091172
$ws.Range("M105").Value = -48089425

$ws.Range("H107").Value = 1171.6666
$ws.Range("I107").Value = 754.61536
$ws.Range("K107").Value = 754.61536
$ws.Range("M107").Value = 1165.38464

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = ""

$ws.Range("H16").Value = 33334330
$ws.Range("I16").Value = 47620016
$ws.Range("K16").Value = 47620016
$ws.Range("M16").Value = -47619729

$ws.Range("H99").Value = 1572.6666
$ws.Range("J99").Value = 1581.3334
$ws.Range("L99").Value = 1581.3334
$ws.Range("N99").Value = -4577.3334

$ws.Range("H113").Value = 33334330
$ws.Range("I113").Value = 47620016
$ws.Range("K113").Value = 47620016
$ws.Range("M113").Value = -47617846

$ws.Range("H126").Value = 1572.6666
$ws.Range("J126").Value = 1581.3334
$ws.Range("L126").Value = 4744.0002
$ws.Range("N126").Value = -9684.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 924.8
$ws.Range("I5").Value = 1028.8148
$ws.Range("J5").Value = 573.75
$ws.Range("K5").Value = 3086.4444
$ws.Range("L5").Value = 1721.25
$ws.Range("M5").Value = -2974.4444
$ws.Range("N5").Value = -1945.25

$ws.Range("H13").Value = 305.6
$ws.Range("I13").Value = 132.25
$ws.Range("J13").Value = 999
$ws.Range("K13").Value = 396.75
$ws.Range("L13").Value = 2997
$ws.Range("M13").Value = -228.75
$ws.Range("N13").Value = -3333

$ws.Range("H131").Value = 13158791
$ws.Range("J131").Value = 933.6197
$ws.Range("L131").Value = 2800.8591
$ws.Range("N131").Value = -12880.8591

$ws.Range("H134").Value = 3375.4546
$ws.Range("I134").Value = 2016.25
$ws.Range("J134").Value = 7000
$ws.Range("K134").Value = 6048.75
$ws.Range("L134").Value = 21000
$ws.Range("M134").Value = -978.75
$ws.Range("N134").Value = -31140

$ws.Range("H135").Value = 924.8
$ws.Range("I135").Value = 1028.8148
$ws.Range("J135").Value = 573.75
$ws.Range("K135").Value = 9259.333200000001
$ws.Range("L135").Value = 5163.75
$ws.Range("M135").Value = -6724.333200000001
$ws.Range("N135").Value = -10233.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2619.8333
$ws.Range("I122").Value = 2832.1177
$ws.Range("J122").Value = 2104.2856
$ws.Range("K122").Value = 8496.3531
$ws.Range("L122").Value = 6312.8568
$ws.Range("M122").Value = -6046.3531
$ws.Range("N122").Value = -11212.8568

$ws.Range("H124").Value = 56000
$ws.Range("J124").Value = 56000
$ws.Range("L124").Value = 56000
$ws.Range("N124").Value = -65820

$ws.Range("H127").Value = 37187.25
$ws.Range("J127").Value = 37187.25
$ws.Range("L127").Value = 37187.25
$ws.Range("N127").Value = -47107.25

$ws.Range("H132").Value = 3220.4443
$ws.Range("I132").Value = 3166.4375
$ws.Range("K132").Value = 9499.3125
$ws.Range("M132").Value = -6969.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1034.4546
$ws.Range("I22").Value = 963.3333
$ws.Range("J22").Value = 1119.8
$ws.Range("K22").Value = 963.3333
$ws.Range("L22").Value = 1119.8
$ws.Range("M22").Value = -668.3333
$ws.Range("N22").Value = -1709.8

$ws.Range("H27").Value = 1034.4546
$ws.Range("I27").Value = 963.3333
$ws.Range("J27").Value = 1119.8
$ws.Range("K27").Value = 963.3333
$ws.Range("L27").Value = 1119.8
$ws.Range("M27").Value = -856.3333
$ws.Range("N27").Value = -1333.8

$ws.Range("H46").Value = 4745
$ws.Range("J46").Value = 5194.4443
$ws.Range("L46").Value = 5194.4443
$ws.Range("N46").Value = -5570.4443

$ws.Range("H61").Value = 1267.5
$ws.Range("I61").Value = 1154.091
$ws.Range("J61").Value = 1683.3334
$ws.Range("K61").Value = 1154.091
$ws.Range("L61").Value = 1683.3334
$ws.Range("M61").Value = -952.0909999999999
$ws.Range("N61").Value = -2087.3334

$ws.Range("H113").Value = 1267.5
$ws.Range("I113").Value = 1154.091
$ws.Range("J113").Value = 1683.3334
$ws.Range("K113").Value = 1154.091
$ws.Range("L113").Value = 1683.3334
$ws.Range("M113").Value = 1015.909
$ws.Range("N113").Value = -6023.3334

$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").Value = ""

$ws.Range("H132").Value = 1712.3704
$ws.Range("I132").Value = 1343.1666
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 4029.4998
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -1499.4998
$ws.Range("N132").Value = -19058

$ws.Range("H133").Value = 35205
$ws.Range("J133").Value = 35205
$ws.Range("L133").Value = 35205
$ws.Range("N133").Value = -45325
